# Apply the edit described by the diff:
# 1. Add a new "Player Info" worksheet as the first sheet with player details
#    (ID / NAME / BATTING_HAND / BOWL_STYLE), styled like the other sheets'
#    header row.
# 2. Rename the "MATCH_CARD_LINK" header to "MATCH_CODE" on the "ODI Batting"
#    (column D) and "ODI Bowling" (column B) sheets.
# 3. Replace the full scorecard URL values in that column with just the
#    numeric match code extracted from the URL.

$wb = $excel.ActiveWorkbook

# --- Build the new "Player Info" sheet -------------------------------------
# Duplicate "ODI Batting" (placed immediately before it) so the new sheet's
# header row ends up sharing the exact same cell style (bold font + border +
# center/top alignment) already used by the other sheets' header rows,
# instead of constructing an equivalent-looking but distinct style.
$battingSheetForInsert = $wb.Worksheets.Item("ODI Batting")
$battingSheetForInsert.Copy($battingSheetForInsert)

$playerInfo = $wb.Worksheets.Item(1)
$playerInfo.Name = "Player Info"

# Drop the copied "ODI Batting" data/columns we don't need.
$playerInfo.Range("E1:J9").Clear()
$playerInfo.Range("A3:D9").Clear()

# Headers (row 1 keeps the inherited header style).
$playerInfo.Range("A1").Value = "ID"
$playerInfo.Range("B1").Value = "NAME"
$playerInfo.Range("C1").Value = "BATTING_HAND"
$playerInfo.Range("D1").Value = "BOWL_STYLE"

# Data row. ID is numeric-looking text ("5774"), so force text storage via a
# temporary Text number format and then drop the format again (mirrors
# typing into a Text-formatted cell in real Excel) so the cell ends up
# string-typed with the plain/default style, matching the other sheets'
# unstyled data cells.
$idCell = $playerInfo.Range("A2")
$idCell.NumberFormat = "@"
$idCell.Value = "5774"
$idCell.ClearFormats()

$playerInfo.Range("B2").Value = "Kyle Alex Jamieson"
$playerInfo.Range("C2").Value = "Right Handed"
$playerInfo.Range("D2").Value = "Right Arm Fast Medium"

$playerInfo.Range("A1").Select() | Out-Null

# --- Re-fetch the other sheets by name --------------------------------------
# (references captured before the sheet-collection change above now point at
# the wrong positional slot, since this host resolves earlier handles by
# index rather than identity)
$battingSheet = $wb.Worksheets.Item("ODI Batting")
$bowlingSheet = $wb.Worksheets.Item("ODI Bowling")

# --- Update "ODI Batting" sheet: column D header + values ------------------
$battingSheet.Range("D1").Value = "MATCH_CODE"
for ($r = 2; $r -le 9; $r++) {
    $cell = $battingSheet.Cells.Item($r, 4)
    $url = $cell.Value2
    if ($url -match "MatchCode=(\d+)") {
        $code = $Matches[1]
        $cell.NumberFormat = "@"
        $cell.Value = $code
        $cell.ClearFormats()
    }
}

# --- Update "ODI Bowling" sheet: column B header + values ------------------
$bowlingSheet.Range("B1").Value = "MATCH_CODE"
for ($r = 2; $r -le 9; $r++) {
    $cell = $bowlingSheet.Cells.Item($r, 2)
    $url = $cell.Value2
    if ($url -match "MatchCode=(\d+)") {
        $code = $Matches[1]
        $cell.NumberFormat = "@"
        $cell.Value = $code
        $cell.ClearFormats() | Out-Null
    }
}
